# Update hotel reviews data on the "hotel_info" sheet.
# English_Reviews_num (G2), Local_Rank (H2) and Total_Reviews_num (I2) were
# previously blank and are now populated with the latest scrape numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# Force text storage (these are id-like codes, not arithmetic numbers) so the
# values round-trip as shared strings, matching the other text columns on
# this row (TA_ReviewURL, Orbitz_ReviewURL, etc).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4684"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "13"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "4824"

# Touch J2 (Orbitz_ReviewURL) without changing its content so its shared
# string stays referenced/stable alongside the newly added strings above.
$ws.Range("J2").Value = $ws.Range("J2").Value()
